$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$meta = $wb.Worksheets.Item(1)

# Update Version value (row 3, column B)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Sheet 2: "Elements" ---
$elem = $wb.Worksheets.Item(2)

# Add constraint text to the InFulfillmentOf.typeId row (row 5), Constraint(s) column (AJ)
$elem.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
